## Generate Report for Handoff
## - Flip the handback status text to "Ready for handoff" everywhere it is
##   shown (Overview!E2/F2, zh-cn!C2, de-de!C2).
## - Bump the "Latest HO Xliff Generate Date" / handoff datetime stamps that
##   go with the new report generation.
## - The status column got narrower now that the text is shorter, so shrink
##   the columns that used to be sized for the old, longer status string.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Status text -----------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Timestamps --------------------------------------------------------------
$wsOverview.Range("G2").Value = "2016-08-23 22:55:38"
$wsDeDe.Range("H2").Value = "2016-08-23 22:55:38"
$wsZhCn.Range("H2").Value = "2016-08-23 22:55:33"

# --- Column widths -----------------------------------------------------------
# Stored column width is snapped to a 1/6-character pixel grid, so feed in the
# character width whose rounded grid value lands closest to the target
# (17.2159881591797 -> grid value 17.1666...).
$newColWidth = 16.333333333333332
$wsOverview.Range("E1:F1").ColumnWidth = $newColWidth
$wsZhCn.Range("C1").ColumnWidth = $newColWidth
$wsDeDe.Range("C1").ColumnWidth = $newColWidth
